$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new trade row (row 9) matching the existing table's data/shape.
$ws.Range("A9").Value = 9597.3799999999992
$ws.Range("B9").Value = 9794.24
$ws.Range("C9").Value = 79.650000000000006
$ws.Range("D9").Value = 78.05
$ws.Range("E9").Value = $false
$ws.Range("F9").Value = -2.0099999999999998
$ws.Range("G9").Value = 42612.67292824074
$ws.Range("H9").Value = $false

# Reuse the date/time number formatting from the row above (G column)
# instead of inventing a new custom format code.
$ws.Range("G8").Copy($ws.Range("G9"))
$ws.Range("G9").Value = 42612.67292824074
